$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking 'cryptos' sheet refresh: updated Price (D) and Volume(1h) (E)
# text columns to the latest scraped snapshot.

# A handful of the new Price strings parse as plain decimals (e.g. '300.93'),
# which Excel's General format would auto-convert to a Number on entry. The
# source data keeps every Price cell as literal text, so pre-format those
# specific cells as Text before writing the value, preserving exact formatting
# (e.g. trailing zeros) and cell type.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated Price (D) / Volume(1h) (E) values
$ws.Range("D2").Value = '43.042.64'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '2.304.40'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '300.93'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = '98.31'
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("D7").Value = '0.518'
$ws.Range("E7").Value = '  +2.35%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").Value = '36.19'
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '17.85'
$ws.Range("E12").Value = '  +1.36%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '2.663.19'
$ws.Range("D16").Value = '2.305.72'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").Value = '43.001.74'
$ws.Range("D19").Value = '12.70'
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = '6.13'
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("D22").Value = '68.30'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").Value = '242.56'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").Value = '166.99'
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").Value = '33.20'
$ws.Range("E32").Value = '  -2.77%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '5.04'
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("D35").Value = '4.75'
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("E39").Value = '  -1.71%  '
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("D43").Value = '2.006.05'
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("D46").Value = '10.19'
$ws.Range("E46").Value = '  +0.81%  '
$ws.Range("D47").Value = '17.59'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  -2.17%  '
$ws.Range("D49").Value = '53.69'
$ws.Range("E49").Value = '  -2.20%  '
$ws.Range("D50").Value = '2.528.97'
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("D51").Value = '72.81'
$ws.Range("E51").Value = '  -5.21%  '
